$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "The "
$ws.Range("B2").Value = "shiny "
$ws.Range("C2").Value = "gondola "
$ws.Range("D2").Value = "moved "
$ws.Range("E2").Value = "slowly."
$ws.Range("F2").Value = "infreq_1"
$ws.Range("G2").Value = "infrequent"
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = "critical"
$ws.Range("J2").Value = 1

# Row 3
$ws.Range("A3").Value = "The "
$ws.Range("B3").Value = "large "
$ws.Range("C3").Value = "mosque "
$ws.Range("D3").Value = "remained "
$ws.Range("E3").Value = "mostly empty."
$ws.Range("F3").Value = "infreq_2"
$ws.Range("G3").Value = "infrequent"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = "critical"
$ws.Range("J3").Value = 2

# Row 4
$ws.Range("A4").Value = "The "
$ws.Range("B4").Value = "noisy "
$ws.Range("C4").Value = "chicken "
$ws.Range("D4").Value = "chased "
$ws.Range("E4").Value = "the sparrows."
$ws.Range("F4").Value = "freq_3"
$ws.Range("G4").Value = "frequent"
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = "critical"
$ws.Range("J4").Value = 3

# Row 5
$ws.Range("A5").Value = "The "
$ws.Range("B5").Value = "sandy "
$ws.Range("C5").Value = "beach "
$ws.Range("D5").Value = "stretched "
$ws.Range("E5").Value = "for many miles."
$ws.Range("F5").Value = "freq_4"
$ws.Range("G5").Value = "frequent"
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = "critical"
$ws.Range("J5").Value = 4

# Update the active selection to C3, matching the saved view state.
$ws.Range("C3").Select()
